$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "วางรากฐาน"

# --- New content below the existing "rule of 3" mission section ---
# (values are written in this exact order so that the shared-strings
#  table is built in the same sequence as the authored workbook)
$ws.Range("B33").Value = "3.feedback ว่า3สิ่งนั้น ใช้พลังงาน สมาธิ เวลา เท่าไหร่"

$ws.Range("A35").Value = "ภารกิจช่วงเวลาทอง"
$ws.Range("B36").Value = "พลังงาน"
$ws.Range("C37").Value = "พยายามงด น้ำตาล คาเฟอีน"
$ws.Range("B38").Value = "เวลา"
$ws.Range("C39").Value = "บันทึก 3 อย่างทุกต้นชม."
$ws.Range("C40").Value = "1.ระดับพลังงาน"
$ws.Range("C41").Value = "2.กำลังทำอะไร"
$ws.Range("C42").Value = "3.อู้งานกี่นาที (โดยประมาณ)"

# Table header row (row 44) - column D/E/F first, then B, then C which
# reuses the already-created "เวลา" string
$ws.Range("D44").Value = "ระดับพลังงาน"
$ws.Range("E44").Value = "กำลังทำอะไร"
$ws.Range("F44").Value = "อู้งานกี่นาที (โดยประมาณ)"
$ws.Range("B44").Value = "วัน"
$ws.Range("C44").Value = "เวลา"

# Style the header row like a "Heading 2" with a thick accent underline
$ws.Range("B44:F44").Style = "Heading 2"
$ws.Rows.Item(44).RowHeight = 18

# Data rows under the header (hour-by-hour log table)
$ws.Range("B45").Value = 21
$ws.Range("C45").Value = 7
$ws.Range("C46").Value = 8
$ws.Range("C47").Value = 9
$ws.Range("C48").Value = 10
$ws.Range("C49").Value = 11
$ws.Range("C50").Value = 12
$ws.Range("C51").Value = 13
$ws.Range("C52").Value = 14
$ws.Range("C53").Value = 15
$ws.Range("C54").Value = 16
$ws.Range("C55").Value = 17
$ws.Range("C56").Value = 18
$ws.Range("C57").Value = 19
$ws.Range("C58").Value = 20
$ws.Range("C59").Value = 21
$ws.Range("C60").Value = 22
$ws.Range("C61").Value = 23

# Column widths for the new table
$ws.Columns.Item(3).ColumnWidth = 9.85546875
$ws.Columns.Item(4).ColumnWidth = 16.42578125
$ws.Columns.Item(5).ColumnWidth = 86
$ws.Columns.Item(6).ColumnWidth = 32.42578125

# Leave the view scrolled/selected near the new table
$ws.Range("D53").Select() | Out-Null
